# Add an "Authors" credit textbox to the title slide (slide 1), placed to
# the right of/below the existing "OpenDSS: Version ..." textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Target position/size (EMU, converted to points -- PowerPoint COM uses points)
$emuPerPoint = 12700.0
$left   = 7130642 / $emuPerPoint
$top    = 6042139 / $emuPerPoint
$width  = 3838167 / $emuPerPoint
$height = 369332  / $emuPerPoint

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 3"

# No wrapping, shape auto-sizes to fit the text (matches a:spAutoFit / wrap="none")
$tb.TextFrame.WordWrap = $false
$tb.TextFrame.AutoSize = 1

# No shape fill
$tb.Fill.Visible = $false

$tb.TextFrame.TextRange.Text = "Authors: Paulo Radatz and Celso Rocha"
